# Applies the dated-worksheet content refresh: new date line plus all 100
# addition/subtraction problems in the table, in an order chosen so that no
# earlier replacement text can be mistakenly matched by a later search
# (e.g. "7+34=" is a substring of "17+34=", so the "17+34=" cell must be
# updated before we search for "7+34=").
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-12 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-13 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("8+65=", $true, $false, $false, $false, $false, $true, 1, $false, "54+1=", 2) | Out-Null
$d.Content.Find.Execute("38+46=", $true, $false, $false, $false, $false, $true, 1, $false, "56+33=", 2) | Out-Null
$d.Content.Find.Execute("95-48=", $true, $false, $false, $false, $false, $true, 1, $false, "78+7=", 2) | Out-Null
$d.Content.Find.Execute("42+23=", $true, $false, $false, $false, $false, $true, 1, $false, "34-1=", 2) | Out-Null
$d.Content.Find.Execute("65+22=", $true, $false, $false, $false, $false, $true, 1, $false, "35+15=", 2) | Out-Null
$d.Content.Find.Execute("15+21=", $true, $false, $false, $false, $false, $true, 1, $false, "85-61=", 2) | Out-Null
$d.Content.Find.Execute("10+5=", $true, $false, $false, $false, $false, $true, 1, $false, "27+43=", 2) | Out-Null
$d.Content.Find.Execute("31-23=", $true, $false, $false, $false, $false, $true, 1, $false, "1+60=", 2) | Out-Null
$d.Content.Find.Execute("54-1=", $true, $false, $false, $false, $false, $true, 1, $false, "30+35=", 2) | Out-Null
$d.Content.Find.Execute("25+35=", $true, $false, $false, $false, $false, $true, 1, $false, "46-28=", 2) | Out-Null
$d.Content.Find.Execute("91-39=", $true, $false, $false, $false, $false, $true, 1, $false, "70-4=", 2) | Out-Null
$d.Content.Find.Execute("8-1=", $true, $false, $false, $false, $false, $true, 1, $false, "68+4=", 2) | Out-Null
$d.Content.Find.Execute("93-13=", $true, $false, $false, $false, $false, $true, 1, $false, "79-3=", 2) | Out-Null
$d.Content.Find.Execute("93-62=", $true, $false, $false, $false, $false, $true, 1, $false, "81-30=", 2) | Out-Null
$d.Content.Find.Execute("6+74=", $true, $false, $false, $false, $false, $true, 1, $false, "59-10=", 2) | Out-Null
$d.Content.Find.Execute("40-7=", $true, $false, $false, $false, $false, $true, 1, $false, "50+6=", 2) | Out-Null
$d.Content.Find.Execute("81-12=", $true, $false, $false, $false, $false, $true, 1, $false, "77-30=", 2) | Out-Null
$d.Content.Find.Execute("33+30=", $true, $false, $false, $false, $false, $true, 1, $false, "66-15=", 2) | Out-Null
$d.Content.Find.Execute("75-17=", $true, $false, $false, $false, $false, $true, 1, $false, "28-3=", 2) | Out-Null
$d.Content.Find.Execute("25+15=", $true, $false, $false, $false, $false, $true, 1, $false, "59-30=", 2) | Out-Null
$d.Content.Find.Execute("1+9=", $true, $false, $false, $false, $false, $true, 1, $false, "16+2=", 2) | Out-Null
$d.Content.Find.Execute("50+31=", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=", 2) | Out-Null
$d.Content.Find.Execute("21+4=", $true, $false, $false, $false, $false, $true, 1, $false, "59+29=", 2) | Out-Null
$d.Content.Find.Execute("60-18=", $true, $false, $false, $false, $false, $true, 1, $false, "86-21=", 2) | Out-Null
$d.Content.Find.Execute("70+20=", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=", 2) | Out-Null
$d.Content.Find.Execute("73-44=", $true, $false, $false, $false, $false, $true, 1, $false, "77+21=", 2) | Out-Null
$d.Content.Find.Execute("17+34=", $true, $false, $false, $false, $false, $true, 1, $false, "87-61=", 2) | Out-Null
$d.Content.Find.Execute("7+34=", $true, $false, $false, $false, $false, $true, 1, $false, "62+2=", 2) | Out-Null
$d.Content.Find.Execute("22+53=", $true, $false, $false, $false, $false, $true, 1, $false, "84-38=", 2) | Out-Null
$d.Content.Find.Execute("47+37=", $true, $false, $false, $false, $false, $true, 1, $false, "45-5=", 2) | Out-Null
$d.Content.Find.Execute("58-38=", $true, $false, $false, $false, $false, $true, 1, $false, "63-14=", 2) | Out-Null
$d.Content.Find.Execute("8+76=", $true, $false, $false, $false, $false, $true, 1, $false, "58-26=", 2) | Out-Null
$d.Content.Find.Execute("61-32=", $true, $false, $false, $false, $false, $true, 1, $false, "86-28=", 2) | Out-Null
$d.Content.Find.Execute("38-11=", $true, $false, $false, $false, $false, $true, 1, $false, "94-3=", 2) | Out-Null
$d.Content.Find.Execute("62-40=", $true, $false, $false, $false, $false, $true, 1, $false, "10+36=", 2) | Out-Null
$d.Content.Find.Execute("28-8=", $true, $false, $false, $false, $false, $true, 1, $false, "6+64=", 2) | Out-Null
$d.Content.Find.Execute("84+12=", $true, $false, $false, $false, $false, $true, 1, $false, "69-26=", 2) | Out-Null
$d.Content.Find.Execute("88-37=", $true, $false, $false, $false, $false, $true, 1, $false, "54-7=", 2) | Out-Null
$d.Content.Find.Execute("37-4=", $true, $false, $false, $false, $false, $true, 1, $false, "82-46=", 2) | Out-Null
$d.Content.Find.Execute("5+91=", $true, $false, $false, $false, $false, $true, 1, $false, "63+6=", 2) | Out-Null
$d.Content.Find.Execute("38-6=", $true, $false, $false, $false, $false, $true, 1, $false, "10+60=", 2) | Out-Null
$d.Content.Find.Execute("98-87=", $true, $false, $false, $false, $false, $true, 1, $false, "45-2=", 2) | Out-Null
$d.Content.Find.Execute("25+25=", $true, $false, $false, $false, $false, $true, 1, $false, "25+37=", 2) | Out-Null
$d.Content.Find.Execute("64-53=", $true, $false, $false, $false, $false, $true, 1, $false, "60-15=", 2) | Out-Null
$d.Content.Find.Execute("21+38=", $true, $false, $false, $false, $false, $true, 1, $false, "77-1=", 2) | Out-Null
$d.Content.Find.Execute("4+57=", $true, $false, $false, $false, $false, $true, 1, $false, "56-5=", 2) | Out-Null
$d.Content.Find.Execute("2+61=", $true, $false, $false, $false, $false, $true, 1, $false, "48+37=", 2) | Out-Null
$d.Content.Find.Execute("51+41=", $true, $false, $false, $false, $false, $true, 1, $false, "87-36=", 2) | Out-Null
$d.Content.Find.Execute("57+6=", $true, $false, $false, $false, $false, $true, 1, $false, "62-2=", 2) | Out-Null
$d.Content.Find.Execute("92+0=", $true, $false, $false, $false, $false, $true, 1, $false, "44-41=", 2) | Out-Null
$d.Content.Find.Execute("66-34=", $true, $false, $false, $false, $false, $true, 1, $false, "16+38=", 2) | Out-Null
$d.Content.Find.Execute("19+64=", $true, $false, $false, $false, $false, $true, 1, $false, "25+68=", 2) | Out-Null
$d.Content.Find.Execute("30-13=", $true, $false, $false, $false, $false, $true, 1, $false, "42-17=", 2) | Out-Null
$d.Content.Find.Execute("28+45=", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=", 2) | Out-Null
$d.Content.Find.Execute("98-69=", $true, $false, $false, $false, $false, $true, 1, $false, "47+25=", 2) | Out-Null
$d.Content.Find.Execute("28+44=", $true, $false, $false, $false, $false, $true, 1, $false, "31-11=", 2) | Out-Null
$d.Content.Find.Execute("42+33=", $true, $false, $false, $false, $false, $true, 1, $false, "33+52=", 2) | Out-Null
$d.Content.Find.Execute("12+72=", $true, $false, $false, $false, $false, $true, 1, $false, "42+29=", 2) | Out-Null
$d.Content.Find.Execute("75-68=", $true, $false, $false, $false, $false, $true, 1, $false, "66+23=", 2) | Out-Null
$d.Content.Find.Execute("94-84=", $true, $false, $false, $false, $false, $true, 1, $false, "43-2=", 2) | Out-Null
$d.Content.Find.Execute("7+1=", $true, $false, $false, $false, $false, $true, 1, $false, "1+37=", 2) | Out-Null
$d.Content.Find.Execute("94-16=", $true, $false, $false, $false, $false, $true, 1, $false, "0+41=", 2) | Out-Null
$d.Content.Find.Execute("91-90=", $true, $false, $false, $false, $false, $true, 1, $false, "83-45=", 2) | Out-Null
$d.Content.Find.Execute("22+0=", $true, $false, $false, $false, $false, $true, 1, $false, "51+15=", 2) | Out-Null
$d.Content.Find.Execute("40-28=", $true, $false, $false, $false, $false, $true, 1, $false, "99-6=", 2) | Out-Null
$d.Content.Find.Execute("42+30=", $true, $false, $false, $false, $false, $true, 1, $false, "86-45=", 2) | Out-Null
$d.Content.Find.Execute("75-31=", $true, $false, $false, $false, $false, $true, 1, $false, "37-10=", 2) | Out-Null
$d.Content.Find.Execute("8+50=", $true, $false, $false, $false, $false, $true, 1, $false, "27+62=", 2) | Out-Null
$d.Content.Find.Execute("63-36=", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=", 2) | Out-Null
$d.Content.Find.Execute("78-51=", $true, $false, $false, $false, $false, $true, 1, $false, "56-14=", 2) | Out-Null
$d.Content.Find.Execute("27-7=", $true, $false, $false, $false, $false, $true, 1, $false, "4+11=", 2) | Out-Null
$d.Content.Find.Execute("62+22=", $true, $false, $false, $false, $false, $true, 1, $false, "27-17=", 2) | Out-Null
$d.Content.Find.Execute("62-30=", $true, $false, $false, $false, $false, $true, 1, $false, "74-27=", 2) | Out-Null
$d.Content.Find.Execute("10+49=", $true, $false, $false, $false, $false, $true, 1, $false, "71-7=", 2) | Out-Null
$d.Content.Find.Execute("97-21=", $true, $false, $false, $false, $false, $true, 1, $false, "98-22=", 2) | Out-Null
$d.Content.Find.Execute("88+7=", $true, $false, $false, $false, $false, $true, 1, $false, "56+37=", 2) | Out-Null
$d.Content.Find.Execute("37-19=", $true, $false, $false, $false, $false, $true, 1, $false, "66+6=", 2) | Out-Null
$d.Content.Find.Execute("98-45=", $true, $false, $false, $false, $false, $true, 1, $false, "46-38=", 2) | Out-Null
$d.Content.Find.Execute("89-45=", $true, $false, $false, $false, $false, $true, 1, $false, "32+53=", 2) | Out-Null
$d.Content.Find.Execute("30+24=", $true, $false, $false, $false, $false, $true, 1, $false, "22-16=", 2) | Out-Null
$d.Content.Find.Execute("18+57=", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=", 2) | Out-Null
$d.Content.Find.Execute("69+19=", $true, $false, $false, $false, $false, $true, 1, $false, "29+48=", 2) | Out-Null
$d.Content.Find.Execute("98-61=", $true, $false, $false, $false, $false, $true, 1, $false, "83+7=", 2) | Out-Null
$d.Content.Find.Execute("11+23=", $true, $false, $false, $false, $false, $true, 1, $false, "18+0=", 2) | Out-Null
$d.Content.Find.Execute("10-0=", $true, $false, $false, $false, $false, $true, 1, $false, "70+23=", 2) | Out-Null
$d.Content.Find.Execute("30-29=", $true, $false, $false, $false, $false, $true, 1, $false, "76+18=", 2) | Out-Null
$d.Content.Find.Execute("42+38=", $true, $false, $false, $false, $false, $true, 1, $false, "17+71=", 2) | Out-Null
$d.Content.Find.Execute("22+14=", $true, $false, $false, $false, $false, $true, 1, $false, "41+26=", 2) | Out-Null
$d.Content.Find.Execute("86-12=", $true, $false, $false, $false, $false, $true, 1, $false, "63-45=", 2) | Out-Null
$d.Content.Find.Execute("42+36=", $true, $false, $false, $false, $false, $true, 1, $false, "25-21=", 2) | Out-Null
$d.Content.Find.Execute("73-70=", $true, $false, $false, $false, $false, $true, 1, $false, "79-26=", 2) | Out-Null
$d.Content.Find.Execute("22-4=", $true, $false, $false, $false, $false, $true, 1, $false, "48-33=", 2) | Out-Null
$d.Content.Find.Execute("7+61=", $true, $false, $false, $false, $false, $true, 1, $false, "99-97=", 2) | Out-Null
$d.Content.Find.Execute("44-29=", $true, $false, $false, $false, $false, $true, 1, $false, "71+21=", 2) | Out-Null
$d.Content.Find.Execute("57-29=", $true, $false, $false, $false, $false, $true, 1, $false, "3+74=", 2) | Out-Null
$d.Content.Find.Execute("53-43=", $true, $false, $false, $false, $false, $true, 1, $false, "81-66=", 2) | Out-Null
$d.Content.Find.Execute("33+40=", $true, $false, $false, $false, $false, $true, 1, $false, "48-5=", 2) | Out-Null
$d.Content.Find.Execute("21+55=", $true, $false, $false, $false, $false, $true, 1, $false, "30+51=", 2) | Out-Null
$d.Content.Find.Execute("49+22=", $true, $false, $false, $false, $false, $true, 1, $false, "36+63=", 2) | Out-Null
$d.Content.Find.Execute("72-41=", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=", 2) | Out-Null
